# Automatic update of files.
# 1) Column C ("Förändrad") on every data row moves from 45184 to 45186
#    (the sheet's "last changed" date serial advances by two days).
# 2) The HYPERLINK() formulas in columns S, T, V, W, X, Y (present only on
#    the rows that actually have attachments) gain a friendly-text second
#    argument equal to the row's case id (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 74
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = $firstRow; $row -le $lastRow; $row++) {

    # -- 1) bump the "Förändrad" date serial in column C --
    $ws.Cells.Item($row, 3).Value2 = 45186

    # -- 2) add the display-text argument to any HYPERLINK formulas on this row --
    $id = $ws.Range("A$row").Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$row")
        $oldFormula = $cell.Formula
        if ($oldFormula -and $oldFormula.Length -gt 0 -and $oldFormula.ToUpper().Contains("HYPERLINK(") -and -not $oldFormula.Contains(",")) {
            $newFormula = $oldFormula.Substring(0, $oldFormula.Length - 1) + ', "' + $id + '")'
            $cell.Formula = $newFormula
        }
    }
}
